$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Miscellaneous")

# Enter the new "Name"/"Specific" values in the same order the author typed
# them (this governs the order new entries land in sharedStrings.xml).
$ws.Range("C4").Value = "Text"
$ws.Range("D4").Value = "insertion"
$ws.Range("D5").Value = "clickEvent"
$ws.Range("D6").Value = "hoverEvent"
$ws.Range("D7").Value = "translate Component"
$ws.Range("D8").Value = "keybind Component"
$ws.Range("D9").Value = "score Component"
$ws.Range("F5").Value = "convert to better names, Text.d.ts"
$ws.Range("F6").Value = "convert to better names, implement item NBT, Tekst.d.ts"
$ws.Range("C11").Value = "CustomError"
$ws.Range("D14").Value = "expectationNotMet"
$ws.Range("D11").Value = "Not implemented"
$ws.Range("D12").Value = "Action can't be performed right now"
$ws.Range("D13").Value = "Can't set"
$ws.Range("D10").Value = "chatToArray("

# Status column (E) values reuse already-existing shared strings.
$ws.Range("E4").Value = "fully implemented"
$ws.Range("E5").Value = "partly implemented"
$ws.Range("E6").Value = "not started"
$ws.Range("E7").Value = "not started"
$ws.Range("E8").Value = "not started"
$ws.Range("E9").Value = "not started"
$ws.Range("E10").Value = "not started"
$ws.Range("E11").Value = "not started"
$ws.Range("E12").Value = "not started"
$ws.Range("E13").Value = "not started"
$ws.Range("E14").Value = "fully implemented"

# Update selected cell on Miscellaneous sheet to match saved selection
$ws.Range("F13").Select()
